$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.138.95'
$ws.Range('E2').Value = '  -1.91%  '
$ws.Range('D3').Value = '1.839.09'
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('E4').Value = '  -0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '239.97'
$c.NumberFormat = 'General'
$ws.Range('E6').Value = '  -2.89%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.06%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2992'
$c.NumberFormat = 'General'
$ws.Range('E8').Value = '  -2.93%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07462'
$c.NumberFormat = 'General'
$ws.Range('E9').Value = '  -4.05%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '23.25'
$c.NumberFormat = 'General'
$ws.Range('E10').Value = '  -2.14%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07647'
$c.NumberFormat = 'General'
$ws.Range('E11').Value = '  -2.37%  '
$ws.Range('D12').Value = '1.836.46'
$ws.Range('E12').Value = '  -1.54%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '5.033'
$c.NumberFormat = 'General'
$ws.Range('E13').Value = '  -2.83%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.6811'
$c.NumberFormat = 'General'
$ws.Range('E14').Value = '  -2.26%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '87.50'
$c.NumberFormat = 'General'
$ws.Range('E15').Value = '  -5.90%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '6.161'
$c.NumberFormat = 'General'
$ws.Range('E16').Value = '  -7.27%  '
$ws.Range('D17').Value = '29.134.61'
$ws.Range('E17').Value = '  -1.89%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000008222'
$c.NumberFormat = 'General'
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('D19').Value = '2.085.53'
$ws.Range('E19').Value = '  -1.30%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '230.60'
$c.NumberFormat = 'General'
$ws.Range('E20').Value = '  -5.54%  '
$ws.Range('E21').Value = '  -2.46%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.NumberFormat = 'General'
$ws.Range('E22').Value = '  -0.03%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '7.345'
$c.NumberFormat = 'General'
$ws.Range('E23').Value = '  -4.07%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.00%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '161.16'
$c.NumberFormat = 'General'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('E26').Value = '  -6.01%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.707'
$c.NumberFormat = 'General'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '18.06'
$c.NumberFormat = 'General'
$ws.Range('E28').Value = '  -1.88%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.505'
$c.NumberFormat = 'General'
$ws.Range('E29').Value = '  -2.83%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.254'
$c.NumberFormat = 'General'
$ws.Range('E30').Value = '  -0.63%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.135'
$c.NumberFormat = 'General'
$ws.Range('E31').Value = '  -1.78%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.195'
$c.NumberFormat = 'General'
$ws.Range('E32').Value = '  -0.38%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.05355'
$c.NumberFormat = 'General'
$ws.Range('E33').Value = '  +4.72%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.7543'
$c.NumberFormat = 'General'
$ws.Range('E34').Value = '  -4.51%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.848'
$c.NumberFormat = 'General'
$ws.Range('E35').Value = '  -3.66%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.132'
$c.NumberFormat = 'General'
$ws.Range('E36').Value = '  -2.60%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.684'
$c.NumberFormat = 'General'
$ws.Range('D38').Value = '1.314.37'
$ws.Range('E38').Value = '  -2.11%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01828'
$c.NumberFormat = 'General'
$ws.Range('E39').Value = '  -3.24%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.720'
$c.NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.84%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.9421'
$c.NumberFormat = 'General'
$ws.Range('E41').Value = '  -2.11%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '6.054'
$c.NumberFormat = 'General'
$ws.Range('E42').Value = '  +0.16%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '105.10'
$c.NumberFormat = 'General'
$ws.Range('E43').Value = '  -1.65%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.9992'
$c.NumberFormat = 'General'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.07919'
$c.NumberFormat = 'General'
$ws.Range('E45').Value = '  +24.29%  '
$ws.Range('D46').Value = '1.985.17'
$ws.Range('E46').Value = '  -1.49%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.5182'
$c.NumberFormat = 'General'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('E48').Value = '  -3.83%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '64.31'
$c.NumberFormat = 'General'
$ws.Range('E49').Value = '  -1.65%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.774'
$c.NumberFormat = 'General'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '9.425'
$c.NumberFormat = 'General'
